# DC - actualización de excel con las rubricas de juego
#
# Adds the data for three new "minor scale" challenges (12, 13, 14) that were
# previously blank placeholders in the rubric sheet, narrows column B a
# little, and leaves the cursor / selection where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Desafio 12 (row 57-61): "Juego de tonos y semitonos en escalas
#     menores..." already had its question text - it was only missing the
#     Pregunta/Respuesta pair.
$ws.Range("C59").Value = "Rem"
$ws.Range("D59").Value = "Re:T,Mi:S,Fa:T,Sol:T,La:S,Sib:T,Do:T"

# --- Desafio 13 (row 62-66): "Reconocimiento auditivo de escalas menores
#     naturales" challenge text plus its Pregunta/Respuesta pair.
$ws.Range("B62").Value = "Reconocimiento auditivo de escalas menores naturales. Escucha una escala menor natural y selecciona su tonalidad correcta."
$ws.Range("C64").Value = "son_esca_do_me"
$ws.Range("D64").Value = "Escala_do_me"

# --- Desafio 14 (row 67-71): "Encuentra la tónica de una escala menor..."
$ws.Range("B67").Value = "Encuentra la tónica de una escala menor. Observa una escala menor y selecciona la nota que actúa como tónica."

# Column B got a touch narrower.
$ws.Columns("B").ColumnWidth = 49.67

# Restore the author's last selection before saving.
$ws.Range("C69").Select()
